# Mark attendance cells as present (1) for the corresponding date rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Invalid + Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count + Real
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count + Real
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Total Attendance Count + Real
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7: Absent
$ws.Range("H7").Value = 1

# Row 8: Absent
$ws.Range("H8").Value = 1

# Row 9: Total Attendance Count + Real
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: Absent
$ws.Range("H10").Value = 1

# Row 11: Absent
$ws.Range("H11").Value = 1

# Row 12: Total Attendance Count + Real
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: Total Attendance Count + Real
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14: Absent
$ws.Range("H14").Value = 1

# Row 15: Absent
$ws.Range("H15").Value = 1

# Row 16: Absent
$ws.Range("H16").Value = 1

# Row 17: Absent
$ws.Range("H17").Value = 1

# Row 18: Absent
$ws.Range("H18").Value = 1
